# Daily attendance processing - 2025-12-08 19:48:15
#
# Normalises the "Recorded By" (column G) audit-trail text: the list of
# recorders for a session is re-ordered (the System/system marker(s) moved
# ahead of the human recorder), matching the source system's new formatting.
# Only cells whose text exactly matches one of the known "before" patterns
# are touched; every other "Recorded By" value (single recorder, or
# combinations such as "admin@admin.com, System") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1

$recordedByCol = 7   # column G = "Recorded By"

$changed = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $current = $cell.Text

    if ($replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
        $changed = $changed + 1
    }
}

Write-Output "Recorded By cells updated: $changed"
